# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout ("K") counts per game row (rows 2-36), replacing the old
# values that had been mistakenly populated from a different stat.
$kValues = @(4, 1, 1, 2, 2, 2, 4, 3, 7, 8, 5, 3, 8, 2, 5, 3, 5, 9, 7, 5, 7, 0, 10, 2, 3, 5, 2, 4, 3, 6, 7, 5, 3, 3, 1)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
